# Refresh the cryptocurrency price/volume table on the active sheet, as
# produced by the nightly "Updated cryptos list" GitHub Actions run.
#
# Sheet layout: A=rank index, B=Coin, C=Link, D=Price, E=Volume(1h).
# Most rows only get fresh D (price) / E (1h volume %) values; a handful
# of rows (33/34, 40/41, 50/51) changed ranking order between runs, so
# their whole row (Coin, Link, Price, Volume) is rewritten.
#
# D-column prices are stored as literal text in the source workbook (e.g.
# "11.50", "1.00", thousand-separator style "66.226.38"), so Set-PriceText
# temporarily forces text format ("@") before writing the new value —
# otherwise Excel would silently coerce these numeric-looking strings into
# real numbers and lose significant trailing digits/formatting — then
# restores the cell's original number format afterwards.
function Set-PriceText($sheet, $addr, $text) {
    $origFmt = $sheet.Range($addr).NumberFormat
    $sheet.Range($addr).NumberFormat = "@"
    $sheet.Range($addr).Value = $text
    $sheet.Range($addr).NumberFormat = $origFmt
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2: update D,E
    Set-PriceText $ws "D2" "66.226.38"
    $ws.Range("E2").Value = "  +0.23%  "
    # Row 3: update D,E
    Set-PriceText $ws "D3" "3.567.29"
    $ws.Range("E3").Value = "  +2.44%  "
    # Row 4: update E
    $ws.Range("E4").Value = "  +0.05%  "
    # Row 5: update D,E
    Set-PriceText $ws "D5" "608.43"
    $ws.Range("E5").Value = "  +0.80%  "
    # Row 6: update D,E
    Set-PriceText $ws "D6" "145.37"
    $ws.Range("E6").Value = "  +1.60%  "
    # Row 7: update D,E
    Set-PriceText $ws "D7" "3.566.30"
    $ws.Range("E7").Value = "  +2.48%  "
    # Row 8: update E
    $ws.Range("E8").Value = "  +0.09%  "
    # Row 9: update D,E
    Set-PriceText $ws "D9" "0.491"
    $ws.Range("E9").Value = "  +3.53%  "
    # Row 10: update E
    $ws.Range("E10").Value = "  +1.15%  "
    # Row 11: update E
    $ws.Range("E11").Value = "  -3.37%  "
    # Row 12: update E
    $ws.Range("E12").Value = "  +0.53%  "
    # Row 13: update D,E
    Set-PriceText $ws "D13" "4.175.54"
    $ws.Range("E13").Value = "  +2.64%  "
    # Row 14: update E
    $ws.Range("E14").Value = "  +2.45%  "
    # Row 15: update E
    $ws.Range("E15").Value = "  -1.12%  "
    # Row 16: update D,E
    Set-PriceText $ws "D16" "3.572.74"
    $ws.Range("E16").Value = "  +2.52%  "
    # Row 17: update D,E
    Set-PriceText $ws "D17" "66.356.30"
    $ws.Range("E17").Value = "  +0.33%  "
    # Row 18: update E
    $ws.Range("E18").Value = "  -1.05%  "
    # Row 19: update D,E
    Set-PriceText $ws "D19" "11.50"
    $ws.Range("E19").Value = "  +10.41%  "
    # Row 20: update E
    $ws.Range("E20").Value = "  +1.04%  "
    # Row 21: update D,E
    Set-PriceText $ws "D21" "14.86"
    $ws.Range("E21").Value = "  +0.77%  "
    # Row 22: update D,E
    Set-PriceText $ws "D22" "429.81"
    $ws.Range("E22").Value = "  +1.97%  "
    # Row 23: update E
    $ws.Range("E23").Value = "  +4.56%  "
    # Row 24: update D,E
    Set-PriceText $ws "D24" "79.20"
    $ws.Range("E24").Value = "  +2.20%  "
    # Row 25: update D,E
    Set-PriceText $ws "D25" "3.712.97"
    $ws.Range("E25").Value = "  +2.82%  "
    # Row 26: update E
    $ws.Range("E26").Value = "  +0.03%  "
    # Row 27: update D,E
    Set-PriceText $ws "D27" "0.0000118"
    $ws.Range("E27").Value = "  +3.77%  "
    # Row 28: update D,E
    Set-PriceText $ws "D28" "2.51"
    $ws.Range("E28").Value = "  +2.18%  "
    # Row 29: update D,E
    Set-PriceText $ws "D29" "7.94"
    $ws.Range("E29").Value = "  -0.12%  "
    # Row 30: update E
    $ws.Range("E30").Value = "  -2.62%  "
    # Row 31: update E
    $ws.Range("E31").Value = "  +0.10%  "
    # Row 32: update D,E
    Set-PriceText $ws "D32" "25.64"
    $ws.Range("E32").Value = "  +1.98%  "
    # Row 33: update B,C,D,E
    $ws.Range("B33").Value = "RenzoRestakedETH"
    $ws.Range("C33").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
    Set-PriceText $ws "D33" "3.566.12"
    $ws.Range("E33").Value = "  +2.57%  "
    # Row 34: update B,C,D,E
    $ws.Range("B34").Value = "Fetch.AI"
    $ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
    Set-PriceText $ws "D34" "1.46"
    $ws.Range("E34").Value = "  -1.54%  "
    # Row 37: update E
    $ws.Range("E37").Value = "  +1.50%  "
    # Row 38: update E
    $ws.Range("E38").Value = "  +2.48%  "
    # Row 39: update D,E
    Set-PriceText $ws "D39" "5.60"
    $ws.Range("E39").Value = "  +0.30%  "
    # Row 40: update B,C,D,E
    $ws.Range("B40").Value = "FirstDigitalUSD"
    $ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    Set-PriceText $ws "D40" "1.00"
    $ws.Range("E40").Value = "  +0.11%  "
    # Row 41: update B,C,D,E
    $ws.Range("B41").Value = "Monero"
    $ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    Set-PriceText $ws "D41" "177.34"
    $ws.Range("E41").Value = "  +3.77%  "
    # Row 42: update D,E
    Set-PriceText $ws "D42" "0.0848"
    $ws.Range("E42").Value = "  -1.78%  "
    # Row 43: update E
    $ws.Range("E43").Value = "  +2.62%  "
    # Row 44: update E
    $ws.Range("E44").Value = "  +0.84%  "
    # Row 45: update D,E
    Set-PriceText $ws "D45" "1.95"
    $ws.Range("E45").Value = "  +1.08%  "
    # Row 46: update D,E
    Set-PriceText $ws "D46" "46.18"
    $ws.Range("E46").Value = "  +2.27%  "
    # Row 47: update E
    $ws.Range("E47").Value = "  +0.75%  "
    # Row 48: update D,E
    Set-PriceText $ws "D48" "25.67"
    $ws.Range("E48").Value = "  -1.97%  "
    # Row 49: update E
    $ws.Range("E49").Value = "  +3.01%  "
    # Row 50: update B,C,D,E
    $ws.Range("B50").Value = "EnergySwap"
    $ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    Set-PriceText $ws "D50" "23.59"
    $ws.Range("E50").Value = "  +9.47%  "
    # Row 51: update B,C,D,E
    $ws.Range("B51").Value = "Cosmos"
    $ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
    Set-PriceText $ws "D51" "7.15"
    $ws.Range("E51").Value = "  +0.43%  "
